$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1476674106001415
$ws.Range("C2").Value = 2.942101459957778
$ws.Range("D2").Value = 0.07049155245767012
$ws.Range("B3").Value = 0.09768040281875964
$ws.Range("C3").Value = 3.298972515789647
$ws.Range("D3").Value = 0.1435797463073023
$ws.Range("B4").Value = 0.1349109666706689
$ws.Range("C4").Value = 2.885949886071024
$ws.Range("D4").Value = 0.06293077365701315
$ws.Range("B5").Value = 0.1176141782755219
$ws.Range("C5").Value = 2.314064943909416
$ws.Range("D5").Value = 0.05760386809086761
$ws.Range("B6").Value = 0.08387920535328325
$ws.Range("C6").Value = 1.891600025215038
$ws.Range("D6").Value = 0.07017953711186031
$ws.Range("B7").Value = 0.06133465385190478
$ws.Range("C7").Value = 1.670071560772693
$ws.Range("D7").Value = 0.1188793223012413
$ws.Range("B8").Value = 0.05568351709503638
$ws.Range("C8").Value = 0.6403596860693295
$ws.Range("D8").Value = 0.10346647721759
$ws.Range("B9").Value = 0.08424191696994403
$ws.Range("C9").Value = 0.07147129295388653
$ws.Range("D9").Value = 0.1325046779212274
$ws.Range("B10").Value = 0.1053977716841419
$ws.Range("C10").Value = -0.4586304841802679
$ws.Range("D10").Value = 0.08549481488210758
$ws.Range("B11").Value = 0.09885937618367405
$ws.Range("C11").Value = -1.214177754731012
$ws.Range("D11").Value = 0.08559483556625029
$ws.Range("B12").Value = 0.06450622686451341
$ws.Range("C12").Value = -2.048689657360445
$ws.Range("D12").Value = 0.1468509823372504
$ws.Range("B13").Value = 0.05095252640635244
$ws.Range("C13").Value = -2.310243552827969
$ws.Range("D13").Value = 0.08806520750585126
$ws.Range("B14").Value = 0.1050488427212418
$ws.Range("C14").Value = -2.727182624186739
$ws.Range("D14").Value = 0.1498751106465459
$ws.Range("B15").Value = 0.1188826385697819
$ws.Range("C15").Value = -2.827916073640093
$ws.Range("D15").Value = 0.1173844685118637
$ws.Range("B16").Value = 0.07541223863442555
$ws.Range("C16").Value = -2.82997844365193
$ws.Range("D16").Value = 0.07962616928736238
$ws.Range("B17").Value = 0.1391262935900381
$ws.Range("C17").Value = -2.955531292268895
$ws.Range("D17").Value = 0.1181665086434596
$ws.Range("B18").Value = 0.09538742837929842
$ws.Range("C18").Value = -2.595698837342453
$ws.Range("D18").Value = 0.06236158478798554
$ws.Range("B19").Value = 0.06376667845470585
$ws.Range("C19").Value = -2.564387514058391
$ws.Range("D19").Value = 0.08550120459922199
$ws.Range("B20").Value = 0.05946586234230136
$ws.Range("C20").Value = -2.22144093092511
$ws.Range("D20").Value = 0.09083280187346386
$ws.Range("B21").Value = 0.14143163881128
$ws.Range("C21").Value = -1.608876037668836
$ws.Range("D21").Value = 0.12823967821486
$ws.Range("B22").Value = 0.08305853201782236
$ws.Range("C22").Value = -1.074418688859741
$ws.Range("D22").Value = 0.08843313942692847
$ws.Range("B23").Value = 0.1244079366159185
$ws.Range("C23").Value = -0.2753698464313605
$ws.Range("D23").Value = 0.05820950321167653
$ws.Range("B24").Value = 0.09447986381113339
$ws.Range("C24").Value = 0.1747408497613457
$ws.Range("D24").Value = 0.1498086926009261
$ws.Range("B25").Value = 0.102706168991809
$ws.Range("C25").Value = 1.28647804516334
$ws.Range("D25").Value = 0.1099705306713902
$ws.Range("B26").Value = 0.08190103872214494
$ws.Range("C26").Value = 1.246365693771726
$ws.Range("D26").Value = 0.1010696732095713
$ws.Range("B27").Value = 0.08040126220466116
$ws.Range("C27").Value = 1.812025271099527
$ws.Range("D27").Value = 0.08616791133993917
$ws.Range("B28").Value = 0.08297304313685699
$ws.Range("C28").Value = 2.034382655514586
$ws.Range("D28").Value = 0.05793545829138692
$ws.Range("B29").Value = 0.1224683169770897
$ws.Range("C29").Value = 2.840469187707176
$ws.Range("D29").Value = 0.1134330744292962
$ws.Range("B30").Value = 0.1401174123704053
$ws.Range("C30").Value = 2.855748465557786
$ws.Range("D30").Value = 0.1338161944184134
$ws.Range("B31").Value = 0.07151188914532082
$ws.Range("C31").Value = 3.098250279686679
$ws.Range("D31").Value = 0.08157090338010761
